# Replace all occurrences of "OIE" with "WOAH" in the relevant text cells
# of the "Sheet 1" and "References" worksheets, per the commit:
#   "OIE replaced with WOAH all Excels"
# Hyperlink target URLs (e.g. https://www.oie.int/...) are left untouched,
# only the visible text content is updated.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet 1")
$refs   = $wb.Worksheets.Item("References")

$nbsp = [char]0x00A0

$sheet1.Range("E5").Value  = "Based on official disease reports to the WOAH"
$sheet1.Range("E6").Value  = "ND is a disease listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code and must be reported to the WOAH. The map to the right displays outbreak points reported to the WOAH early warning system since 2005."
$sheet1.Range("E7").Value  = "As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:"
$sheet1.Range("E14").Value = "Countries are coloured according to the available information regarding their stable disease situation (disease status legend). This information is provided by countries through the WOAH monitoring system, which is a different reporting channel.<br>Immediate notifications (points) and disease status (country/region colours) are reported to the WOAH in different spatial and temporal scales, and therefore are displayed in the map as layers which can be filtered independently."
$sheet1.Range("E17").Value = "For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}."
$sheet1.Range("E21").Value = "A summary of the disease in animal hosts is given in the {ref008:WOAH Technical disease card}."
$sheet1.Range("E34").Value = "Humans may become infected and get conjunctivitis, but the condition is generally very mild and self-limiting ({ref010:WOAH}). Visit the Public Health impact section of this story map."
$sheet1.Range("E43").Value = "Refer to the {ref008:WOAH Technical disease card} for a key summary of the virus characteristics. "
$sheet1.Range("E55").Value = "Refer to the {ref008:WOAH Technical disease card} for a key summary of the disease transmission and epidemiological parameters."
$sheet1.Range("E67").Value = "WOAH-prescribed tests for detection of the agent include: virus isolation, conventional PCR and real-time-PCR. Assessment of the virus virulence is required  and this shall be based on gene sequencing or the intracerebral pathogenicity index (ICPI). WOAH-prescribed tests for detection of immune response include a variety of commercial enzyme-linked immunosorbent assays (ELISA), and haemagglutination inhibition (HI) tests  ({ref010:WOAH," + $nbsp + "Terrestrial Manual})"
$sheet1.Range("E87").Value = "Click on the figure to the right to visit the original source of this information (WOAH)."
$sheet1.Range("E88").Value = "ND_WOAH.html"
$sheet1.Range("E98").Value = "Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data."

$refs.Range("C2").Value  = "WOAH-WAHIS (WOAH World Animal Health Information System)"
$refs.Range("C6").Value  = "WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$refs.Range("C9").Value  = "WOAH (World Organisation for Animal Health) Technical Disease Card: Newcastle Disease. 2013."
$refs.Range("C10").Value = "WOAH (World Organisation for Animal Health), 2021. Newcastle Disease. Chapter 10.9. WOAH Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$refs.Range("C11").Value = "WOAH (World Organisation for Animal Health), 2019. Newcastle Disease. Chapter 3.3.14 WOAH Terrestrial Manual 2021. WOAH, Paris, France"
